$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A10 and A11 content
$a10 = $ws.Range("A10").Value2
$a11 = $ws.Range("A11").Value2
$ws.Range("A10").Value = $a11
$ws.Range("A11").Value = $a10

# Row 15: replace error message text
$ws.Range("A15").Value = "Foutmelding wordt weergegeven op de lasergun door middel van LED"

# Row 34: replace HP/sound text
$ws.Range("A34").Value = "Wanneer een speler geen hitpoints meer heeft (en dus 'dood' is) wordt dit bekend gemaakt door een geluidssignaal van 1000Hz en kan de speler het wapen 3 seconden niet gebruiken"

# New column C width (bestFit-style autofit approximation)
$ws.Columns("C").ColumnWidth = 21.17

# Update selection to match final cursor position
$ws.Range("C14").Select()

Write-Output "done"
